# Remove the stray blank line / page-break / footer paragraph block that
# used to sit between the "Requisitos" list and the next page's content:
#
#   ...LOB1036: Geometria Analítica (Requisito fraco)
#   <blank paragraph>                                   <- remove
#   <blank paragraph, page-break-before>                <- remove
#   © 2020 . Contact: luizeleno@usp.br. ...              <- remove
#   <blank paragraph>
#   <blank paragraph, page-break-before>
#
$d = $word.ActiveDocument

# Locate the anchor paragraph ("LOB1036: ...") so the deletion is tied to
# content rather than to a fragile hard-coded paragraph index.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1036: Geometria Anal*tica (Requisito fraco)*") {
        $anchor = $p
        break
    }
}
if ($anchor -eq $null) {
    throw "Could not find the 'LOB1036: Geometria Analítica' paragraph"
}

# The three paragraphs to remove are the three immediately following the
# anchor paragraph: a blank paragraph, a blank page-break paragraph, and
# the copyright/footer paragraph.
$firstToDelete = $anchor.Next()
$secondToDelete = $firstToDelete.Next()
$thirdToDelete = $secondToDelete.Next()

if ($thirdToDelete.Range.Text -notlike "*Contact: luizeleno@usp.br*") {
    throw "Unexpected document structure - footer paragraph not found where expected"
}

$delRange = $d.Range($firstToDelete.Range.Start, $thirdToDelete.Range.End)
$delRange.Delete()
